# NSMB - Begin 8-8
# Appends the newly-tracked checkpoints/splits for level 8-8 (rows 175-191)
# to the "V4" results sheet, extending the running IF(B>0,C-B,0) delta
# formula in column D down through the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 175; A = "Checkpoint 791";      B = 56265; C = 66317 },
    @{ Row = 176; A = "Land on Koopa";       B = 56378; C = 66431 },
    @{ Row = 177; A = "Checkpoint 1154";     B = 56428; C = 66483 },
    @{ Row = 178; A = "Checkpoint 1702";     B = 56640; C = 66695 },
    @{ Row = 179; A = "Checkpoint 1927";     B = 56715; C = 66770 },
    @{ Row = 180; A = "Checkpoint 2230";     B = 56898; C = 66954 },
    @{ Row = 181; A = "Checkpoint 2421";     B = 56982; C = 67039 },
    @{ Row = 182; A = "Checkpoint 2550";     B = 57020; C = 67079 },
    @{ Row = 183; A = "Checkpoint ";         B = 57085; C = 67144 },
    @{ Row = 184; A = "Checkpoint 2927";     B = 57134; C = 67193 },
    @{ Row = 185; A = "Checkpoint 3271";     B = 57234; C = 67293 },
    @{ Row = 186; A = "Checkpoint 3576";     B = 57309; C = 67368 },
    @{ Row = 187; A = "Blast out of pipe";   B = 57531; C = 67591 },
    @{ Row = 188; A = "Get Flag";            B = 57771; C = 67831 },
    @{ Row = 189; A = "End lLevel";          B = 58289; C = 68349 },
    @{ Row = 190; A = "Enter 8-8";           B = 58622; C = 69018 },
    @{ Row = 191; A = "1st Move";            B = 58848; C = 69265 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = $r.A
    $ws.Range("B$i").Value = $r.B
    $ws.Range("C$i").Value = $r.C
    $ws.Range("D$i").Formula = "=IF(B$i>0,C$i-B$i,0)"
}

# Restore the frozen header pane / scroll the view down to the newly
# added rows and leave the selection where the author left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 175
$ws.Range("B192").Select()
